$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Regional association study" heading: demote from Heading1 to
#    Heading2 (word-format tweak mentioned in the commit message).
#    The paragraph already carries the "regional-association-study"
#    bookmark, so use it to locate the paragraph robustly.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks.Item("regional-association-study")
$bm.Range.Paragraphs.Item(1).Style = "Heading 2"

# ---------------------------------------------------------------------
# 2) Fix the typo: "1-Mb region" -> "10-Mb region".
# ---------------------------------------------------------------------
$d.Content.Find.Execute("1-Mb region", $true, $false, $false, $false, $false,
                         $true, 1, $false, "10-Mb region", 2)

# ---------------------------------------------------------------------
# 3) Insert a new "Figure legend" Heading1 paragraph (with a matching
#    bookmark) right before the final "Figure: ..." caption paragraph.
# ---------------------------------------------------------------------
$figIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Figure: Regional association scanning results")) {
        $figIndex = $i
        break
    }
}

$figPara = $d.Paragraphs.Item($figIndex)
$figPara.Range.InsertParagraphBefore()

# The freshly inserted (still empty) paragraph now occupies the slot
# that used to hold the "Figure: ..." paragraph.
$newPara = $d.Paragraphs.Item($figIndex)
$newPara.Range.Text = "Figure legend"
$newPara.Style = "Heading 1"

$startPos = $newPara.Range.Start
$collapsedRange = $d.Range($startPos, $startPos)
$d.Bookmarks.Add("figure-legend", $collapsedRange)
